# Add a "product_code" column (G) to the checklist sheet and populate
# the single known product code (G-2-1) for the Fortune Quest RPG row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in G1
$ws.Range("G1").Value = "product_code"

# Format the data rows of the new column as text (numFmtId 49 / "@")
$ws.Range("G2:G10").NumberFormat = "@"

# Only row 4 (Fortune Quest RPG, 1997) currently has a known product code
$ws.Range("G4").Value = "G-2-1"

# Match the author's final selection in the saved workbook
$ws.Range("G5").Select()
